# Edit: "before using angular-bootstrap ui"
#
# 1) Fix the cypher text for the "date" merge clause so it quotes the
#    `today` placeholder the same way the other clauses do:
#       merge (date:date {date:today})  ->  merge (date:date {date:'"+today+"'})
#    (cell A21 on the "cypher" sheet; its shared-string value feeds the
#    concatenation formula in D21 via =B21&A21&C21).
#
# 2) Row 23 ("create (note)-[linkTo:linkTo]->(keyword) ") was still using
#    a bare opening quote in column B (i.e. it behaved like the first
#    line of a new statement). Switch it to the continuation prefix
#    (+") like every other row in the block, so D23 gets a leading "+".
#
# 3) Update the sheet's active cell/selection to D21 (was A25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cypher")

# --- 1: fix the "today" placeholder quoting ---
$ws.Cells.Item(21, 1).Value = "merge (date:date {date:'`"+today+`"'}) "

# --- 2: row 23, column B -> continuation prefix `+"` (with quote-prefix
# formatting, like the other `+"` cells in the sheet) ---
$ws.Cells.Item(23, 2).Value = "'+`""

# --- 3: move the selection/active cell to D21 ---
$ws.Range("D21").Select()
